# fix: checking date if already exist
# The "15/1/2021" attendance column header was a duplicate/incorrect date;
# it should actually be "14/1/2021". Correct the header and re-derive the
# attendance (P/A) for that date's column for the rows whose record was
# filed under the wrong date, then recompute each row's Total / Percentage
# from the (now corrected) three attendance columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the mislabeled date header in G1.
$ws.Range("G1").Value2 = "14/1/2021"

# 2) Rows whose "15/1/2021"(now "14/1/2021") attendance entry needs to flip.
$rowsToFlip = @(3, 5, 6, 7, 8, 11, 16, 18, 21, 25)

foreach ($r in $rowsToFlip) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $current = $gCell.Value2
    if ($current -eq "P") {
        $gCell.Value2 = "A"
    } else {
        $gCell.Value2 = "P"
    }

    # Recompute Total (C) and Percentage (D) from E, F, G attendance marks.
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2

    $total = 0
    if ($eVal -eq "P") { $total = $total + 1 }
    if ($fVal -eq "P") { $total = $total + 1 }
    if ($gVal -eq "P") { $total = $total + 1 }

    $percentage = ($total * 100) / 3

    $ws.Cells.Item($r, 3).Value2 = $total
    $ws.Cells.Item($r, 4).Value2 = $percentage
}
